$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestForOrderProvisioned")
$ws.Activate()

# Remove the "password" (col B) and "purchesOrderNo" (col C) columns
# entirely, shifting what's left of the sheet down to a single column.
# This also drops the now-unused "password" shared string on save.
$ws.Columns("B:C").Delete()

# The remaining column A becomes the purchase-order test data:
# header in row 1, value in row 2.
$ws.Range("A1").Value = "purchesOrderNo"
$ws.Range("A2").Value = "isp-abc-fgtd-dcrt-tfdtc-asdef"

# Move the active selection to F8, matching the saved workbook state.
$ws.Range("F8").Select()
